$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column retains its text formatting so values such as
# "1.006" or "0.01969" are not reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.169.90"
$ws.Range("E2").Value = "  -1.80%  "
$ws.Range("D3").Value = "1.823.06"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.59%  "
$ws.Range("D5").Value = "312.76"
$ws.Range("E5").Value = "  -2.00%  "
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("D7").Value = "0.4233"
$ws.Range("E7").Value = "  -1.87%  "
$ws.Range("D8").Value = "0.3694"
$ws.Range("E8").Value = "  -1.50%  "
$ws.Range("D9").Value = "0.07240"
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("D10").Value = "0.8594"
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("D11").Value = "20.97"
$ws.Range("E11").Value = "  -2.94%  "
$ws.Range("D12").Value = "1.824.02"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").Value = "6.711"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").Value = "0.07086"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "5.305"
$ws.Range("E15").Value = "  -2.78%  "
$ws.Range("D16").Value = "89.55"
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").Value = "0.000008848"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "15.03"
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("D21").Value = "27.352.25"
$ws.Range("E21").Value = "  -1.19%  "
$ws.Range("D22").Value = "5.132"
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("D23").Value = "10.91"
$ws.Range("E23").Value = "  -2.37%  "
$ws.Range("D24").Value = "2.063.91"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").Value = "1.989"
$ws.Range("D26").Value = "152.61"
$ws.Range("E26").Value = "  -1.90%  "
$ws.Range("D27").Value = "2.200"
$ws.Range("E27").Value = "  +2.91%  "
$ws.Range("D28").Value = "18.44"
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("D29").Value = "5.227"
$ws.Range("E29").Value = "  -3.13%  "
$ws.Range("D30").Value = "116.26"
$ws.Range("E30").Value = "  -3.14%  "
$ws.Range("D31").Value = "0.08837"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").Value = "1.190"
$ws.Range("E32").Value = "  -3.33%  "
$ws.Range("D33").Value = "0.7515"
$ws.Range("E33").Value = "  -3.32%  "
$ws.Range("D34").Value = "4.439"
$ws.Range("E34").Value = "  -2.63%  "
$ws.Range("E37").Value = "  -2.21%  "
$ws.Range("D38").Value = "0.01969"
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("D39").Value = "0.05250"
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("D40").Value = "7.329"
$ws.Range("E40").Value = "  +2.09%  "
$ws.Range("D41").Value = "2.880"
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("D42").Value = "0.1694"
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("D43").Value = "0.5035"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("D44").Value = "8.693"
$ws.Range("E44").Value = "  -2.90%  "
$ws.Range("D45").Value = "10.58"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").Value = "106.68"
$ws.Range("E46").Value = "  -3.57%  "
$ws.Range("D47").Value = "0.4738"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").Value = "1.005"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").Value = "0.06391"
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("D50").Value = "1.667"
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("D51").Value = "1.863"
$ws.Range("E51").Value = "  -1.24%  "

# Row 35/36: Frax and HuobiToken swap identity and values
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.805"
$ws.Range("E35").Value = "  -3.40%  "

$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "1.005"
$ws.Range("E36").Value = "  -0.60%  "
